# [Kadastro App] Yeni kayit eklendi: 3015
# Adds a new record row (row 74) with Kayit No 3015 to both the
# "Kayitlar" master list and the per-birim "Erdemli" sheet.

$wb = $excel.ActiveWorkbook

# Kayit No, Tarih, Birim, Parsel Sayisi, Is, Personeller
$newRecord = @("3015", "2025-09-11", "Erdemli", "1", "3B", "SERDAR ARSLAN (Tekniker), ÖZKAN AKBAŞ (Mühendis)")

$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # -4121 is the xlDown direction constant: walk from the header row to the
    # last contiguous populated row in column A, then append right after it.
    $lastRow = $ws.Cells.Item(1, 1).End(-4121).Row
    $newRowIndex = $lastRow + 1

    # Kayit No, Tarih and Parsel Sayisi look numeric/date-like; force them to
    # stay plain text (matching the rest of the column) with a quote-prefix,
    # just like typing '3015 / '2025-09-11 / '1 into the cell.
    $ws.Cells.Item($newRowIndex, 1).Value = "'" + $newRecord[0]
    $ws.Cells.Item($newRowIndex, 2).Value = "'" + $newRecord[1]
    $ws.Cells.Item($newRowIndex, 3).Value = $newRecord[2]
    $ws.Cells.Item($newRowIndex, 4).Value = "'" + $newRecord[3]
    $ws.Cells.Item($newRowIndex, 5).Value = $newRecord[4]
    $ws.Cells.Item($newRowIndex, 6).Value = $newRecord[5]
}
